$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 7: The Bleat Is On / Maple Wand
$ws.Range("H7").Value = 50000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 50000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 50000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -50224

# Row 14: Wand-full Tonight / Budding Maple Wand
$ws.Range("H14").Value = 50000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 50000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 50000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -50382

# Row 38: Just Give Him a Serum / Hi-Potion of Strength
$ws.Range("H38").Value = 402.66666
$ws.Range("I38").Value = 104
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 312
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = 60
$ws.Range("N38").Value = -3744

# Row 43: Growing Is Knowing / Growth Formula Gamma
$ws.Range("H43").Value = 61405156
$ws.Range("I43").Value = 111112340
$ws.Range("J43").Value = 16668693
$ws.Range("K43").Value = 111112340
$ws.Range("L43").Value = 16668693
$ws.Range("M43").Value = -111112271
$ws.Range("N43").Value = -16668831

# Row 58: A Matter of Vital Importance / Mega-Potion of Vitality
$ws.Range("H58").Value = 3146.45
$ws.Range("I58").Value = 327.41666
$ws.Range("J58").Value = 7375
$ws.Range("K58").Value = 982.2499799999999
$ws.Range("L58").Value = 22125
$ws.Range("M58").Value = -832.2499799999999
$ws.Range("N58").Value = -22425

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 3574432.5
$ws.Range("I132").Value = 3762297.2
$ws.Range("J132").Value = 5003
$ws.Range("K132").Value = 11286891.6
$ws.Range("L132").Value = 15009
$ws.Range("M132").Value = -11284361.6

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 4167761.8
$ws.Range("I137").Value = 996
$ws.Range("J137").Value = 8334527.5
$ws.Range("K137").Value = 2988
$ws.Range("L137").Value = 25003582.5
$ws.Range("M137").Value = -438
$ws.Range("N137").Value = -25008682.5

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 6300
$ws.Range("I138").Value = 6300
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 18900
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -13760
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 15655.787
$ws.Range("I32").Value = 16896.736
$ws.Range("J32").Value = 4487.25
$ws.Range("K32").Value = 16896.736
$ws.Range("L32").Value = 4487.25
$ws.Range("M32").Value = -16609.736

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 4489.8223
$ws.Range("I132").Value = 4342.6875
$ws.Range("J132").Value = 4852
$ws.Range("K132").Value = 13028.0625
$ws.Range("L132").Value = 14556
$ws.Range("M132").Value = -10498.0625
$ws.Range("N132").Value = -19616

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 22435.164
$ws.Range("I134").Value = 30473.143
$ws.Range("J134").Value = 2340.2144
$ws.Range("K134").Value = 91419.429
$ws.Range("L134").Value = 7020.6432
$ws.Range("M134").Value = -88884.429
$ws.Range("N134").Value = -12090.6432

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 4349976.5
$ws.Range("I31").Value = 2147.5715
$ws.Range("J31").Value = 11113266
$ws.Range("K31").Value = 2147.5715
$ws.Range("L31").Value = 11113266
$ws.Range("M31").Value = -1852.5715
$ws.Range("N31").Value = -11113856

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 4349976.5
$ws.Range("I34").Value = 2147.5715
$ws.Range("J34").Value = 11113266
$ws.Range("K34").Value = 2147.5715
$ws.Range("L34").Value = 11113266
$ws.Range("M34").Value = -1945.5715
$ws.Range("N34").Value = -11113670

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 683.92
$ws.Range("I58").Value = 652
$ws.Range("J58").Value = 829.3333
$ws.Range("K58").Value = 652
$ws.Range("L58").Value = 829.3333
$ws.Range("M58").Value = -449
$ws.Range("N58").Value = -1235.3333

# Row 127: In Rod We Trust / Red Pine Fishing Rod
$ws.Range("H127").Value = 33972.832
$ws.Range("I127").Value = 8354.5
$ws.Range("J127").Value = 39096.5
$ws.Range("K127").Value = 8354.5
$ws.Range("L127").Value = 39096.5
$ws.Range("M127").Value = -3394.5
$ws.Range("N127").Value = -49016.5

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 683.92
$ws.Range("I136").Value = 652
$ws.Range("J136").Value = 829.3333
$ws.Range("K136").Value = 1956
$ws.Range("L136").Value = 2487.9999
$ws.Range("M136").Value = 594
$ws.Range("N136").Value = -7587.9999

$ws = $wb.Worksheets.Item("CUL")
# Row 14: Keep Your Powder Dry / Kukuru Powder
$ws.Range("H14").Value = 232.2
$ws.Range("I14").Value = 232.2
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 696.5999999999999
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -523.5999999999999

# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 965
$ws.Range("I122").Value = 966.8889
$ws.Range("J122").Value = 963.3
$ws.Range("K122").Value = 8702.000100000001
$ws.Range("L122").Value = 8669.699999999999
$ws.Range("M122").Value = -6252.000100000001
$ws.Range("N122").Value = -13569.7

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 3006633.5
$ws.Range("I131").Value = 6784.4116
$ws.Range("J131").Value = 5556505.5
$ws.Range("K131").Value = 20353.2348
$ws.Range("L131").Value = 16669516.5
$ws.Range("M131").Value = -15313.2348
$ws.Range("N131").Value = -16679596.5

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 5827.2
$ws.Range("I122").Value = 10066.667
$ws.Range("J122").Value = 3000.889
$ws.Range("K122").Value = 30200.001
$ws.Range("L122").Value = 9002.667000000001
$ws.Range("M122").Value = -27750.001
$ws.Range("N122").Value = -13902.667

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 48624.51
$ws.Range("I132").Value = 75608.41
$ws.Range("J132").Value = 3089.1875
$ws.Range("K132").Value = 226825.23
$ws.Range("L132").Value = 9267.5625
$ws.Range("M132").Value = -224295.23

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 8953.299999999999
$ws.Range("I132").Value = 13200.294
$ws.Range("J132").Value = 3399.5386
$ws.Range("K132").Value = 39600.882
$ws.Range("L132").Value = 10198.6158
$ws.Range("M132").Value = -37070.882

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 5473.75
$ws.Range("I136").Value = 6251.6665
$ws.Range("J136").Value = 3140
$ws.Range("K136").Value = 18754.9995
$ws.Range("L136").Value = 9420
$ws.Range("M136").Value = -16204.9995
$ws.Range("N136").Value = -14520

$ws = $wb.Worksheets.Item("WVR")
# Row 74: Clothing the Naked Truth / Ramie Robe of Casting
$ws.Range("H74").Value = 12152
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 12152
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 12152
$ws.Range("N74").Value = -14024

# Row 77: When in Robes (L) / Ramie Robe of Casting
$ws.Range("H77").Value = 12152
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 12152
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 36456
$ws.Range("N77").Value = -45816

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1003.4909
$ws.Range("I132").Value = 936.55554
$ws.Range("J132").Value = 1130.3158
$ws.Range("K132").Value = 2809.66662
$ws.Range("L132").Value = 3390.9474
$ws.Range("M132").Value = -279.66662
$ws.Range("N132").Value = -8450.947400000001
